# Juno: check in to OLPRODLOC.
#
# Market-trends table corrections on Sheet1 (Table1):
#   - Column header wording tweak: "...の売上..." -> "...の販売..." for the
#     "Artisanal Chai" column.
#   - April (row 5), "事前に作成された..." column: was mistakenly entered as
#     the clock-time text "4:36" -- fix it back to the numeric 436.
#   - May (row 6): "Artisanal Chai" and "事前に作成された..." columns were
#     retyped as clock-time text "3:15" / "5:48".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Restore D5 to a real number while the workbook is still in automatic
# calculation mode, so the SUM() formula in B5 (which depends on D5)
# recomputes correctly.
$ws.Range("D5").Value = 436

# Switch to manual calculation before touching C6/D6 below: those two become
# text, which the sheet's "=SUM(Cn+Dn)" helper column can't add, and the
# source data keeps the row's previously-computed total untouched rather
# than flipping it to an error.
$excel.Calculation = -4135

# Header: "売上" (sales) -> "販売" (selling) for the Artisanal Chai column.
$ws.Range("C1").Value = "Artisanal Chai の販売 (ユニット数)"

# May row: values re-entered as time-like text.
$ws.Range("C6").Value = "3:15"
$ws.Range("D6").Value = "5:48"
